# Append two new job rows (inserted right after the existing header/first
# data row), refresh the "fetched at" timestamp for every data row, and
# rewire the hyperlinks for column F so they line up with the new layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTimestamp = "2025-09-15 01:19:14"

# Remove every existing hyperlink first so we can rebuild the collection
# cleanly (ids re-number sequentially as we re-add them below).
$ws.Hyperlinks.Delete()

# Shift the existing rows 3-5 down to rows 5-7, opening up rows 3-4 for the
# two newly scraped listings.
$ws.Rows("3:4").Insert()

# --- Row 2 : existing listing, timestamp refreshed only ---------------
$ws.Cells.Item(2, 1).Value = $newTimestamp

# --- Row 3 (new) : 屋上貸切露天風呂の空き状況確認システム開発 ---------
$ws.Cells.Item(3, 1).Value = $newTimestamp
$ws.Cells.Item(3, 2).Value = "【急募】屋上貸切露天風呂の空き状況確認システム開発"
$ws.Cells.Item(3, 3).Value = "システム開発"
$ws.Cells.Item(3, 4).Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Cells.Item(3, 5).Value = "期限情報なし"
$ws.Cells.Item(3, 6).Value = "https://www.lancers.jp/work/detail/5389645"
$ws.Cells.Item(3, 7).Value = 125
$ws.Cells.Item(3, 8).Value = "◆開発,システム開発"

# --- Row 4 (new) : FBA商品(在庫過多商品)をヤフオクで併売するツール開発依頼 ---
$ws.Cells.Item(4, 1).Value = $newTimestamp
$ws.Cells.Item(4, 2).Value = "FBA商品(在庫過多商品)をヤフオクで併売するツール開発依頼"
$ws.Cells.Item(4, 3).Value = "システム開発"
$ws.Cells.Item(4, 4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(4, 5).Value = "期限情報なし"
$ws.Cells.Item(4, 6).Value = "https://www.lancers.jp/work/detail/5393539"
$ws.Cells.Item(4, 7).Value = 123
$ws.Cells.Item(4, 8).Value = "◆ツール,開発"

# --- Rows 5-7 : pre-existing listings (shifted down), timestamp refreshed ---
$ws.Cells.Item(5, 1).Value = $newTimestamp
$ws.Cells.Item(6, 1).Value = $newTimestamp
$ws.Cells.Item(7, 1).Value = $newTimestamp

# --- Rebuild hyperlinks for column F, in row order (F2..F7) -----------
$ws.Hyperlinks.Add($ws.Cells.Item(2, 6), "https://www.lancers.jp/work/detail/5393508")
$ws.Hyperlinks.Add($ws.Cells.Item(3, 6), "https://www.lancers.jp/work/detail/5389645")
$ws.Hyperlinks.Add($ws.Cells.Item(4, 6), "https://www.lancers.jp/work/detail/5393539")
$ws.Hyperlinks.Add($ws.Cells.Item(5, 6), "https://www.lancers.jp/work/detail/5393175")
$ws.Hyperlinks.Add($ws.Cells.Item(6, 6), "https://www.lancers.jp/work/detail/5393406")
$ws.Hyperlinks.Add($ws.Cells.Item(7, 6), "https://www.lancers.jp/work/detail/5393471")

# Re-apply the workbook's "Hyperlink" cell style so column F keeps reusing
# the original style slot instead of the ad-hoc one Hyperlinks.Add creates.
$ws.Cells.Item(2, 6).Style = "Hyperlink"
$ws.Cells.Item(3, 6).Style = "Hyperlink"
$ws.Cells.Item(4, 6).Style = "Hyperlink"
$ws.Cells.Item(5, 6).Style = "Hyperlink"
$ws.Cells.Item(6, 6).Style = "Hyperlink"
$ws.Cells.Item(7, 6).Style = "Hyperlink"
